$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 157:226 originally carry an empty, default-styled B cell (<c r="Bxxx" s="0"/>).
# The edit removes those empty placeholder cells entirely.
$ws.Range("B157:B226").ClearContents()

function Add-WordRow($r, $a, $b, $c) {
    # Inserting the row copies formatting down from the row above (gives col A
    # the wrapped/"s=3" look used throughout the vocabulary table).
    $ws.Rows.Item($r).Insert()
    # Columns B/C should end up with the plain default look ("s=1"), so strip
    # whatever got copied down before typing the new values in.
    $ws.Range("B$r").Clear()
    $ws.Range("C$r").Clear()
    $ws.Range("A$r").Value = $a
    $ws.Range("B$r").Value = $b
    $ws.Range("C$r").Value = $c
}

Add-WordRow 227 "من" "man " "I "
Add-WordRow 228 "تو" "to " "you "
Add-WordRow 229 "او" "oo" "s/he"
Add-WordRow 230 "این/اون" "een/oon" "this/that"
Add-WordRow 231 "ما" "man " "many"
Add-WordRow 232 "شما" "shoma" "shoma"
Add-WordRow 233 "آنها" "anha" "they"

$ws.Range("B233").Select()
